# Generate Report for Handoff
# Replace the old report GUID/hash tokens and timestamps with the new ones.

$wb = $excel.ActiveWorkbook

$oldGuid = "6bfa942e-663c-46bf-b98e-0cac86dc3ee1"
$newGuid = "4598626f-8c43-4320-953b-1f437e950371"

$oldHash = "d8261609b10a70db4e91eb58b2fade30624746b8"
$newHash = "b358715d64d5f018982568cefc120a1f1e855a33"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: A2 (source file name), B2 (hyperlink display path), G2 (latest HO xliff generate date)
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-27 22:56:54"

# zh-cn sheet: A2 (source file name), G2 (latest handoff file), H2 (latest handoff datetime)
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-27 22:56:50"

# de-de sheet: A2 (source file name), G2 (latest handoff file), H2 (latest handoff datetime)
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-27 22:56:54"

